$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition) - first worksheet
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 1101
$ws1.Range("F3").Value = 805
$ws1.Range("F4").Value = 271
$ws1.Range("F5").Value = 45
$ws1.Range("F8").Value = 2038
$ws1.Range("F9").Value = 7573
$ws1.Range("F10").Value = 901
$ws1.Range("F11").Value = 414
$ws1.Range("F12").Value = 345
$ws1.Range("F14").Value = 401
$ws1.Range("F15").Value = 151
$ws1.Range("F16").Value = 7718
$ws1.Range("F17").Value = 304
$ws1.Range("F18").Value = 1340
$ws1.Range("F22").Value = 145
$ws1.Range("F23").Value = 300
$ws1.Range("F24").Value = 137
$ws1.Range("F26").Value = 16
$ws1.Range("F27").Value = 106
$ws1.Range("F29").Value = 409
$ws1.Range("F30").Value = 610
$ws1.Range("F32").Value = 90
$ws1.Range("F33").Value = 59
$ws1.Range("F34").Value = 78
$ws1.Range("F37").Value = 68

# Sheet 4: 全部类型 (All types) - fourth worksheet
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 1101
$ws4.Range("F3").Value = 805
$ws4.Range("F4").Value = 271
$ws4.Range("F5").Value = 45
$ws4.Range("F8").Value = 2039
$ws4.Range("F9").Value = 7573
$ws4.Range("F10").Value = 901
$ws4.Range("F11").Value = 414
$ws4.Range("F12").Value = 345
$ws4.Range("F14").Value = 401
$ws4.Range("F15").Value = 151
$ws4.Range("F16").Value = 7718
$ws4.Range("F17").Value = 304
$ws4.Range("F18").Value = 1340
$ws4.Range("F22").Value = 145
$ws4.Range("F23").Value = 300
$ws4.Range("F24").Value = 137
$ws4.Range("F26").Value = 16
$ws4.Range("F27").Value = 106
$ws4.Range("F29").Value = 409
$ws4.Range("F30").Value = 610
$ws4.Range("F32").Value = 90
$ws4.Range("F33").Value = 59
$ws4.Range("F34").Value = 78
$ws4.Range("F37").Value = 68
